$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.0388964089197
$ws.Range("C2").Value = 5.403934544651637
$ws.Range("D2").Value = 15.05347885465967
$ws.Range("E2").Value = 16.47815336791211
$ws.Range("G2").Value = 3.688789165428493
$ws.Range("J2").Value = 9.414973690434911
$ws.Range("K2").Value = 12.50214419229085
$ws.Range("N2").Value = 20.54798110000768
$ws.Range("O2").Value = 29.41185983877596
$ws.Range("B3").Value = 12.81976026308304
$ws.Range("C3").Value = 5.238687378505033
$ws.Range("D3").Value = 14.99378611203812
$ws.Range("E3").Value = 16.41881176756571
$ws.Range("G3").Value = 3.691032338964267
$ws.Range("J3").Value = 9.423133020389288
$ws.Range("K3").Value = 12.3581250483019
$ws.Range("N3").Value = 20.61290211705236
$ws.Range("O3").Value = 29.47022467434093
$ws.Range("B4").Value = 12.68631091763931
$ws.Range("C4").Value = 5.135544169415604
$ws.Range("D4").Value = 14.96043446935483
$ws.Range("E4").Value = 16.38601898433948
$ws.Range("G4").Value = 3.692482708844998
$ws.Range("J4").Value = 9.429560351096331
$ws.Range("K4").Value = 12.27144249917762
$ws.Range("N4").Value = 20.65461200133004
$ws.Range("O4").Value = 29.51190205620972
$ws.Range("B5").Value = 12.63228010592196
$ws.Range("C5").Value = 5.093157151883576
$ws.Range("D5").Value = 14.9476833672066
$ws.Range("E5").Value = 16.37358143771291
$ws.Range("G5").Value = 3.693092175999174
$ws.Range("J5").Value = 9.432536147212476
$ws.Range("K5").Value = 12.2365969266192
$ws.Range("N5").Value = 20.67207518612818
$ws.Range("O5").Value = 29.53035192033441
$ws.Range("B6").Value = 12.62333168775105
$ws.Range("C6").Value = 5.086099417470838
$ws.Range("D6").Value = 14.9456170642981
$ws.Range("E6").Value = 16.37157236794157
$ws.Range("G6").Value = 3.693194492429137
$ws.Range("J6").Value = 9.433051817594038
$ws.Range("K6").Value = 12.2308408625739
$ws.Range("N6").Value = 20.67500312358938
$ws.Range("O6").Value = 29.53350396130492
$ws.Range("B7").Value = 12.68558072066952
$ws.Range("C7").Value = 5.134973869593671
$ws.Range("D7").Value = 14.96025908991257
$ws.Range("E7").Value = 16.38584748684674
$ws.Range("G7").Value = 3.692490853631954
$ws.Range("J7").Value = 9.429599039666412
$ws.Range("K7").Value = 12.27097057187835
$ws.Range("N7").Value = 20.65484562708847
$ws.Range("O7").Value = 29.51214494573246
$ws.Range("B8").Value = 12.96315402824819
$ws.Range("C8").Value = 5.347345103725434
$ws.Range("D8").Value = 15.03221830937626
$ws.Range("E8").Value = 16.45694236902307
$ws.Range("G8").Value = 3.689547485881892
$ws.Range("J8").Value = 9.417492955882162
$ws.Range("K8").Value = 12.45214786814057
$ws.Range("N8").Value = 20.56998310301365
$ws.Range("O8").Value = 29.43077009161536
$ws.Range("B9").Value = 13.5129861641098
$ws.Range("C9").Value = 5.747753224178631
$ws.Range("D9").Value = 15.19901955894092
$ws.Range("E9").Value = 16.62479213252428
$ws.Range("G9").Value = 3.684352447365993
$ws.Range("J9").Value = 9.404990297428686
$ws.Range("K9").Value = 12.81952739221019
$ws.Range("N9").Value = 20.41816941724725
$ws.Range("O9").Value = 29.31765779476827
$ws.Range("B10").Value = 13.91606944158619
$ws.Range("C10").Value = 6.028861550109064
$ws.Range("D10").Value = 15.33650175827753
$ws.Range("E10").Value = 16.76474044703171
$ws.Range("G10").Value = 3.680883509484763
$ws.Range("J10").Value = 9.402639090406312
$ws.Range("K10").Value = 13.09442544838378
$ws.Range("N10").Value = 20.31544416813565
$ws.Range("O10").Value = 29.26302588563649
$ws.Range("B11").Value = 14.09839038665995
$ws.Range("C11").Value = 6.153311395499269
$ws.Range("D11").Value = 15.40211415186125
$ws.Range("E11").Value = 16.83184650885634
$ws.Range("G11").Value = 3.679380117630209
$ws.Range("J11").Value = 9.403048445766492
$ws.Range("K11").Value = 13.2200580795825
$ws.Range("N11").Value = 20.27060626156308
$ws.Range("O11").Value = 29.2443801361135
$ws.Range("B12").Value = 14.16721208341321
$ws.Range("C12").Value = 6.199901307682707
$ws.Range("D12").Value = 15.42738544734517
$ws.Range("E12").Value = 16.85773704987025
$ws.Range("G12").Value = 3.678821494189046
$ws.Range("J12").Value = 9.403415523662655
$ws.Range("K12").Value = 13.26767328403185
$ws.Range("N12").Value = 20.25389802933725
$ws.Range("O12").Value = 29.23821346537879
$ws.Range("B13").Value = 14.15240092220057
$ws.Range("C13").Value = 6.189891845677827
$ws.Range("D13").Value = 15.42192417370581
$ws.Range("E13").Value = 16.85214002394497
$ws.Range("G13").Value = 3.678941329695399
$ws.Range("J13").Value = 9.403327045875015
$ws.Range("K13").Value = 13.25741735086508
$ws.Range("N13").Value = 20.25748441634317
$ws.Range("O13").Value = 29.2395017829113
$ws.Range("B14").Value = 14.10405716958712
$ws.Range("C14").Value = 6.157155362434084
$ws.Range("D14").Value = 15.40418480885507
$ws.Range("E14").Value = 16.83396703619144
$ws.Range("G14").Value = 3.679333945648304
$ws.Range("J14").Value = 9.40307439919853
$ws.Range("K14").Value = 13.22397483360215
$ws.Range("N14").Value = 20.26922624265888
$ws.Range("O14").Value = 29.24385487413901
$ws.Range("B15").Value = 14.07441466488782
$ws.Range("C15").Value = 6.137032230558239
$ws.Range("D15").Value = 15.39337382623182
$ws.Range("E15").Value = 16.82289743400545
$ws.Range("G15").Value = 3.679575823263161
$ws.Range("J15").Value = 9.402947243285418
$ws.Range("K15").Value = 13.20349439528282
$ws.Range("N15").Value = 20.27645369525113
$ws.Range("O15").Value = 29.24663774323857
$ws.Range("B16").Value = 13.90412735348505
$ws.Range("C16").Value = 6.020655649409236
$ws.Range("D16").Value = 15.33227425491362
$ws.Range("E16").Value = 16.76042282629344
$ws.Range("G16").Value = 3.680983256980906
$ws.Range("J16").Value = 9.402642048358176
$ws.Range("K16").Value = 13.08622307975158
$ws.Range("N16").Value = 20.31841240597484
$ws.Range("O16").Value = 29.26436946365145
$ws.Range("B17").Value = 13.79934411755806
$ws.Range("C17").Value = 5.94835193758567
$ws.Range("D17").Value = 15.29556685987328
$ws.Range("E17").Value = 16.72296726121249
$ws.Range("G17").Value = 3.681865750840576
$ws.Range("J17").Value = 9.402833209819777
$ws.Range("K17").Value = 13.01439901589517
$ws.Range("N17").Value = 20.34463652505795
$ws.Range("O17").Value = 29.27683818079619
$ws.Range("B18").Value = 13.73898129862018
$ws.Range("C18").Value = 5.906444426847563
$ws.Range("D18").Value = 15.27474382611529
$ws.Range("E18").Value = 16.70174889283334
$ws.Range("G18").Value = 3.682380367234198
$ws.Range("J18").Value = 9.403082406803771
$ws.Range("K18").Value = 12.97314482986403
$ws.Range("N18").Value = 20.35989813041802
$ws.Range("O18").Value = 29.28459397720859
$ws.Range("B19").Value = 13.7185294512712
$ws.Range("C19").Value = 5.892201652620644
$ws.Range("D19").Value = 15.26774381200776
$ws.Range("E19").Value = 16.69462102914822
$ws.Range("G19").Value = 3.682555816484606
$ws.Range("J19").Value = 9.403190713173052
$ws.Range("K19").Value = 12.95918798898168
$ws.Range("N19").Value = 20.36509608823854
$ws.Range("O19").Value = 29.28732022634934
$ws.Range("B20").Value = 13.81050873965326
$ws.Range("C20").Value = 5.956082274604518
$ws.Range("D20").Value = 15.29944451733073
$ws.Range("E20").Value = 16.7269209396035
$ws.Range("G20").Value = 3.681771080762958
$ws.Range("J20").Value = 9.402798453009527
$ws.Range("K20").Value = 13.02203920936089
$ws.Range("N20").Value = 20.34182649000482
$ws.Range("O20").Value = 29.27545039781536
$ws.Range("B21").Value = 14.11826338955882
$ws.Range("C21").Value = 6.166785747829598
$ws.Range("D21").Value = 15.4093838813155
$ws.Range("E21").Value = 16.83929201795972
$ws.Range("G21").Value = 3.679218335457145
$ws.Range("J21").Value = 9.403142857571796
$ws.Range("K21").Value = 13.23379693849243
$ws.Range("N21").Value = 20.26577004076125
$ws.Range("O21").Value = 29.24255199011388
$ws.Range("B22").Value = 14.31809075410891
$ws.Range("C22").Value = 6.301348195969967
$ws.Range("D22").Value = 15.48370654354338
$ws.Range("E22").Value = 16.91551653402834
$ws.Range("G22").Value = 3.677612185409945
$ws.Range("J22").Value = 9.404603671214721
$ws.Range("K22").Value = 13.37241179079348
$ws.Range("N22").Value = 20.21764127464197
$ws.Range("O22").Value = 29.22626278561224
$ws.Range("B23").Value = 14.21158088884889
$ws.Range("C23").Value = 6.229830648636881
$ws.Range("D23").Value = 15.44381862672211
$ws.Range("E23").Value = 16.87458492886904
$ws.Range("G23").Value = 3.678463743323444
$ws.Range("J23").Value = 9.403711162007014
$ws.Range("K23").Value = 13.29842426939172
$ws.Range("N23").Value = 20.24318446287247
$ws.Range("O23").Value = 29.23447933070577
$ws.Range("B24").Value = 13.80546158725133
$ws.Range("C24").Value = 5.952588442232189
$ws.Range("D24").Value = 15.29769055333078
$ws.Range("E24").Value = 16.72513249816641
$ws.Range("G24").Value = 3.681813858488905
$ws.Range("J24").Value = 9.402813732625894
$ws.Range("K24").Value = 13.01858495493079
$ws.Range("N24").Value = 20.34309633042793
$ws.Range("O24").Value = 29.27607598495673
$ws.Range("B25").Value = 13.3640960519329
$ws.Range("C25").Value = 5.641509347252727
$ws.Range("D25").Value = 15.15121646178993
$ws.Range("E25").Value = 16.57640901111815
$ws.Range("G25").Value = 3.685696480113529
$ws.Range("J25").Value = 9.407170885203021
$ws.Range("K25").Value = 12.71909287869947
$ws.Range("N25").Value = 20.45768511723422
$ws.Range("O25").Value = 29.53350396130492
